$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Row 7 = "Experimental" property; the Value column (B7) was empty and now must hold
# the literal text "true" (a shared string, not a native boolean, per the FHIR export
# convention used throughout this sheet).
$ws.Range("B7").Value = "'true"

# Re-paste using the formatting of the row above so the quote-prefix style that Value
# just introduced on B7 collapses back to the sheet's normal body style (same as every
# other cell in the column), while the cell keeps holding literal text "true".
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Row 8 = "Date" property; refresh the generation timestamp.
$ws.Range("B8").Value = "2023-02-01T09:05:11-06:00"
